$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add a new (empty) cell O3, matching the formatting of N3
# (bottom-medium border, Times New Roman 9pt, no fill) ---
$ws.Range("N3").Copy() | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null

# --- Row 4: add new year-header cell O4 = 2021, matching N4's formatting ---
$ws.Range("N4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null
$ws.Range("O4").Value = 2021

# --- Row 5: update data values and add new O5 cell ---
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1

$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$ws.Range("O5").Value = 4.0999999999999996

$excel.CutCopyMode = 0

# --- Selection moves to P4 ---
$ws.Range("P4").Select() | Out-Null
